$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.391.39"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.426.29"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'502.22"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "'127.41"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.546"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "2.436.79"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("E11").Value = "  -5.05%  "
$ws.Range("D12").Value = "'5.15"
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("E13").Value = "  -4.83%  "
$ws.Range("D14").Value = "2.859.30"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "57.338.49"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'21.58"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("E17").Value = "  -3.68%  "
$ws.Range("D18").Value = "2.436.92"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("D20").Value = "'312.39"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "'4.06"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'5.67"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "'62.97"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "'0.402"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "'7.13"
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("D29").Value = "'169.11"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").Value = "  -4.69%  "
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "'17.58"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "'36.41"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "'1.43"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "'0.745"
$ws.Range("E41").Value = "  -6.19%  "
$ws.Range("D42").Value = "'268.23"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").Value = "'4.81"
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("D45").Value = "'0.576"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "'118.30"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("D48").Value = "'0.0481"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "'16.99"
$ws.Range("E49").Value = "  -5.04%  "
$ws.Range("E50").Value = "  -3.96%  "
$ws.Range("D51").Value = "'16.40"
$ws.Range("E51").Value = "  -4.53%  "
